$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.4115256666666666
$ws.Range("H2").Value = 1.234577
$ws.Range("I2").Value = 0.2245998342667577
$ws.Range("J2").Value = 0.2245998342667577
$ws.Range("M2").Value = 30.61919899999999
$ws.Range("N2").Value = 91.85759699999998
$ws.Range("O2").Value = 0.6951390881735714
$ws.Range("P2").Value = 0.6951390881735714
$ws.Range("Q2").Value = 12.60058628127433
$ws.Range("R2").Value = 113.405276531469
$ws.Range("S2").Value = 0.1561281239961292
$ws.Range("T2").Value = 0.1561281239961292
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.4115256666666666
$ws.Range("H3").Value = 1.234577
$ws.Range("I3").Value = 0.2245998342667577
$ws.Range("J3").Value = 0.2245998342667577
$ws.Range("O3").Value = 0.09263417906992544
$ws.Range("P3").Value = 0.09263417906992545
$ws.Range("Q3").Value = 1.679153116008
$ws.Range("R3").Value = 15.112378044072
$ws.Range("S3").Value = 0.02080562126654241
$ws.Range("T3").Value = 0.02080562126654241
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.4115256666666666
$ws.Range("H4").Value = 1.234577
$ws.Range("I4").Value = 0.2245998342667577
$ws.Range("J4").Value = 0.2245998342667577
$ws.Range("M4").Value = 9.348075333333332
$ws.Range("N4").Value = 28.04422599999999
$ws.Range("O4").Value = 0.2122267327565031
$ws.Range("P4").Value = 0.2122267327565031
$ws.Range("Q4").Value = 3.846972933600221
$ws.Range("R4").Value = 34.62275640240199
$ws.Range("S4").Value = 0.04766608900408609
$ws.Range("T4").Value = 0.04766608900408609
$ws.Range("I5").Value = 0.3944722233087159
$ws.Range("J5").Value = 0.3944722233087159
$ws.Range("M5").Value = 30.61919899999999
$ws.Range("N5").Value = 91.85759699999998
$ws.Range("O5").Value = 0.6951390881735714
$ws.Range("P5").Value = 0.6951390881735714
$ws.Range("Q5").Value = 22.13083238282366
$ws.Range("R5").Value = 199.177491445413
$ws.Range("S5").Value = 0.2742130616206222
$ws.Range("T5").Value = 0.2742130616206222
$ws.Range("I6").Value = 0.3944722233087159
$ws.Range("J6").Value = 0.3944722233087159
$ws.Range("O6").Value = 0.09263417906992544
$ws.Range("P6").Value = 0.09263417906992545
$ws.Range("S6").Value = 0.0365416105720912
$ws.Range("T6").Value = 0.03654161057209121
$ws.Range("I7").Value = 0.3944722233087159
$ws.Range("J7").Value = 0.3944722233087159
$ws.Range("M7").Value = 9.348075333333332
$ws.Range("N7").Value = 28.04422599999999
$ws.Range("O7").Value = 0.2122267327565031
$ws.Range("P7").Value = 0.2122267327565031
$ws.Range("Q7").Value = 6.756567613150444
$ws.Range("R7").Value = 60.80910851835399
$ws.Range("S7").Value = 0.08371755111600247
$ws.Range("T7").Value = 0.08371755111600247
$ws.Range("G8").Value = 0.6979596666666668
$ws.Range("H8").Value = 2.093879
$ws.Range("I8").Value = 0.3809279424245264
$ws.Range("J8").Value = 0.3809279424245264
$ws.Range("M8").Value = 30.61919899999999
$ws.Range("N8").Value = 91.85759699999998
$ws.Range("O8").Value = 0.6951390881735714
$ws.Range("P8").Value = 0.6951390881735714
$ws.Range("Q8").Value = 21.37096592764033
$ws.Range("R8").Value = 192.338693348763
$ws.Range("S8").Value = 0.26479790255682
$ws.Range("T8").Value = 0.26479790255682
$ws.Range("G9").Value = 0.6979596666666668
$ws.Range("H9").Value = 2.093879
$ws.Range("I9").Value = 0.3809279424245264
$ws.Range("J9").Value = 0.3809279424245264
$ws.Range("O9").Value = 0.09263417906992544
$ws.Range("P9").Value = 0.09263417906992545
$ws.Range("Q9").Value = 2.847893203416
$ws.Range("R9").Value = 25.63103883074401
$ws.Range("S9").Value = 0.03528694723129183
$ws.Range("T9").Value = 0.03528694723129183
$ws.Range("G10").Value = 0.6979596666666668
$ws.Range("H10").Value = 2.093879
$ws.Range("I10").Value = 0.3809279424245264
$ws.Range("J10").Value = 0.3809279424245264
$ws.Range("M10").Value = 9.348075333333332
$ws.Range("N10").Value = 28.04422599999999
$ws.Range("O10").Value = 0.2122267327565031
$ws.Range("P10").Value = 0.2122267327565031
$ws.Range("Q10").Value = 6.524579543628223
$ws.Range("R10").Value = 58.72121589265399
$ws.Range("S10").Value = 0.08084309263641458
$ws.Range("T10").Value = 0.08084309263641458
